$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Budget")
$ws.Activate()

# Update the Year value in A2 from 2023 to 2022
$ws.Range("A2").Value = 2022

# Update the active cell selection to E5
$ws.Range("E5").Select()
